$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "REX_DEF"
for ($r = 2; $r -le 24; $r++) {
    $ws.Range("F$r").Value = "[]"
}
